# civic update + NCI thesaurus 26.01d

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("compounds")
$ws2 = $wb.Worksheets.Item("biomarkers")

# "compounds" sheet, row 3 = NCI Thesaurus source entry.
# Bump its source_version (column E) from "25.12e" to "26.01d".
$ws1.Range("E3").Value = "26.01d"

# View/selection state changes:
# "compounds" becomes the selected/active sheet (it was "biomarkers").
$ws1.Activate()
$ws1.Range("E3").Select() | Out-Null

# "biomarkers" sheet (CIViC row) is no longer the active tab; its
# selection moves from E3 to E6.
$ws2.Range("E6").Select() | Out-Null

# Leave "compounds" as the active sheet when the workbook is saved.
$ws1.Activate() | Out-Null
